# Auto-generated edit script: adds 2022-11-01 violent-crime counts
# across the citywide summary, the by-neighborhood summary, and every
# individual neighborhood sheet that recorded an incident that day.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 6144
$ws.Range("I3").Value = 6397
$ws.Range("C4").Value = 1815
$ws.Range("E4").Value = 1969
$ws.Range("G4").Value = 1445
$ws.Range("I4").Value = 1470
$ws.Range("I6").Value = 7259
$ws.Range("C7").Value = 28358
$ws.Range("E7").Value = 25973
$ws.Range("G7").Value = 24669
$ws.Range("I7").Value = 21862

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 170
$ws.Range("I4").Value = 89
$ws.Range("I7").Value = 689
$ws.Range("I8").Value = 1311
$ws.Range("I9").Value = 110
$ws.Range("I11").Value = 326
$ws.Range("I15").Value = 251
$ws.Range("I19").Value = 604
$ws.Range("I20").Value = 552
$ws.Range("I22").Value = 60
$ws.Range("I27").Value = 192
$ws.Range("I29").Value = 1345
$ws.Range("I31").Value = 219
$ws.Range("I33").Value = 989
$ws.Range("I37").Value = 692
$ws.Range("I42").Value = 767
$ws.Range("I43").Value = 187
$ws.Range("E46").Value = 56
$ws.Range("I47").Value = 156
$ws.Range("I48").Value = 290
$ws.Range("I50").Value = 108
$ws.Range("I51").Value = 255
$ws.Range("I52").Value = 471
$ws.Range("I53").Value = 231
$ws.Range("I54").Value = 447
$ws.Range("C63").Value = 247
$ws.Range("G63").Value = 208
$ws.Range("I63").Value = 70
$ws.Range("I65").Value = 514
$ws.Range("I66").Value = 63
$ws.Range("I67").Value = 843
$ws.Range("I75").Value = 71
$ws.Range("I76").Value = 315
$ws.Range("I77").Value = 138
$ws.Range("I78").Value = 295
$ws.Range("I79").Value = 618
$ws.Range("I83").Value = 478
$ws.Range("I84").Value = 191
$ws.Range("I85").Value = 993
$ws.Range("I90").Value = 275
$ws.Range("I92").Value = 61
$ws.Range("I95").Value = 335
$ws.Range("I96").Value = 238
$ws.Range("C101").Value = 28358
$ws.Range("E101").Value = 25973
$ws.Range("G101").Value = 24669
$ws.Range("I101").Value = 21862

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 284
$ws.Range("I3").Value = 383
$ws.Range("I6").Value = 248
$ws.Range("I7").Value = 993

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I2").Value = 127
$ws.Range("I7").Value = 471

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I2").Value = 135
$ws.Range("I6").Value = 87
$ws.Range("I7").Value = 326

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 394
$ws.Range("I3").Value = 369
$ws.Range("I7").Value = 1311

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I2").Value = 52
$ws.Range("I7").Value = 231

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 226
$ws.Range("I3").Value = 215
$ws.Range("I6").Value = 181
$ws.Range("I7").Value = 689

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I2").Value = 71
$ws.Range("I3").Value = 56
$ws.Range("I7").Value = 238

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I3").Value = 231
$ws.Range("I7").Value = 692

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 310
$ws.Range("I6").Value = 260
$ws.Range("I7").Value = 843

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I4").Value = 12
$ws.Range("I7").Value = 219

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I3").Value = 64
$ws.Range("I7").Value = 191

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 170
$ws.Range("I3").Value = 154
$ws.Range("I4").Value = 21
$ws.Range("I7").Value = 514

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I3").Value = 175
$ws.Range("I7").Value = 478

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I2").Value = 114
$ws.Range("I7").Value = 335

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 223
$ws.Range("I4").Value = 44
$ws.Range("I6").Value = 312
$ws.Range("I7").Value = 989

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I6").Value = 213
$ws.Range("I7").Value = 447

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 395
$ws.Range("I6").Value = 373
$ws.Range("I7").Value = 1345

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 204
$ws.Range("I6").Value = 182
$ws.Range("I7").Value = 604

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I2").Value = 45
$ws.Range("I7").Value = 290

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I2").Value = 64
$ws.Range("I7").Value = 315

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 191
$ws.Range("I3").Value = 237
$ws.Range("I6").Value = 260
$ws.Range("I7").Value = 767

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I3").Value = 73
$ws.Range("I7").Value = 295

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("E4").Value = 7
$ws.Range("E7").Value = 56

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 180
$ws.Range("I6").Value = 183
$ws.Range("I7").Value = 618

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I2").Value = 150
$ws.Range("I7").Value = 552

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I6").Value = 53
$ws.Range("I7").Value = 156

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I6").Value = 95
$ws.Range("I7").Value = 251

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("I2").Value = 31
$ws.Range("I7").Value = 108

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("I6").Value = 26
$ws.Range("I7").Value = 63

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("I3").Value = 38
$ws.Range("I7").Value = 110

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I2").Value = 61
$ws.Range("I7").Value = 170

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("I4").Value = 5
$ws.Range("I7").Value = 61

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("I2").Value = 53
$ws.Range("I7").Value = 192

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("I6").Value = 20
$ws.Range("I7").Value = 71

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I2").Value = 91
$ws.Range("I3").Value = 70
$ws.Range("I7").Value = 275

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I2").Value = 53
$ws.Range("I3").Value = 68
$ws.Range("I7").Value = 255

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I3").Value = 32
$ws.Range("I7").Value = 187

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("I2").Value = 24
$ws.Range("I7").Value = 60

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("I2").Value = 44
$ws.Range("I7").Value = 138

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("I2").Value = 34
$ws.Range("I7").Value = 89
